$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row containing the "超電磁ロボ　コン・バトラーV" post (row 54).
# All subsequent rows shift up by one automatically.
$ws.Rows(54).Delete()
